$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alpha1F")

# Row 16 mirrors the layout of row 15 (next HKL index, reusing the
# "HexGrid-60degTilt5degRes" label) with the averaged intensities set to 1.

# Copy the formatting (bold font + border) from A15 onto A16, then set its value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1
